$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.187431931495667
$ws.Range("B1").Value = 1.960836887359619
$ws.Range("C1").Value = 6.496143817901611
$ws.Range("D1").Value = 2.296064615249634
$ws.Range("E1").Value = 1.195733547210693
